$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Code Review 2" scores (column C) for the four team members
$ws.Range("C2").Value = 25
$ws.Range("C3").Value = 25
$ws.Range("C4").Value = 25
$ws.Range("C5").Value = 25

# Update the active selection
$ws.Range("E10").Select()
